$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C ("Förändrad") rows 2-9 from 45212 to 45221
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 3).Value = 45221
}
